$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.07523809325411
$ws.Range("C2").Value = 9.391514786110049
$ws.Range("E2").Value = 23.2621734548212
$ws.Range("F2").Value = 35.99670634690045
$ws.Range("G2").Value = 18.26812227189786
$ws.Range("H2").Value = 11.72971330077861
$ws.Range("I2").Value = 16.306225043786
$ws.Range("J2").Value = 7.284418104089199
$ws.Range("M2").Value = 19.31117669893643
$ws.Range("O2").Value = 16.43601249144006
$ws.Range("B3").Value = 11.35212310554846
$ws.Range("C3").Value = 8.98383837190249
$ws.Range("E3").Value = 23.2338489291328
$ws.Range("F3").Value = 36.003228662548
$ws.Range("G3").Value = 18.3960433355379
$ws.Range("H3").Value = 11.79055398353849
$ws.Range("I3").Value = 16.44589654945521
$ws.Range("J3").Value = 7.298631964090283
$ws.Range("M3").Value = 19.02447727249596
$ws.Range("O3").Value = 16.54337753169289
$ws.Range("B4").Value = 10.88203695381998
$ws.Range("C4").Value = 8.72294774832355
$ws.Range("E4").Value = 23.22123631705509
$ws.Range("F4").Value = 36.0180941823734
$ws.Range("G4").Value = 18.484990529547
$ws.Range("H4").Value = 11.83028257433835
$ws.Range("I4").Value = 16.53598957653838
$ws.Range("J4").Value = 7.307884849594987
$ws.Range("M4").Value = 18.84817603203729
$ws.Range("O4").Value = 16.61424485473718
$ws.Range("B5").Value = 10.68396523316405
$ws.Range("C5").Value = 8.614084944861343
$ws.Range("E5").Value = 23.21730187198964
$ws.Range("F5").Value = 36.02687917133792
$ws.Range("G5").Value = 18.52382546822396
$ws.Range("H5").Value = 11.84706883689647
$ws.Range("I5").Value = 16.57379572467995
$ws.Range("J5").Value = 7.311787960612288
$ws.Range("M5").Value = 18.77634164368327
$ws.Range("O5").Value = 16.64436352273948
$ws.Range("B6").Value = 10.65068437144067
$ws.Range("C6").Value = 8.595858064082448
$ws.Range("E6").Value = 23.2167214467881
$ws.Range("F6").Value = 36.02850247482139
$ws.Range("G6").Value = 18.53042948112695
$ws.Range("H6").Value = 11.84989221419986
$ws.Range("I6").Value = 16.58013946981514
$ws.Range("J6").Value = 7.312444080990117
$ws.Range("M6").Value = 18.76441644802074
$ws.Range("O6").Value = 16.64943946900512
$ws.Range("B7").Value = 10.87939195813005
$ws.Range("C7").Value = 8.721489742833036
$ws.Range("E7").Value = 23.2211783714091
$ws.Range("F7").Value = 36.01820162421076
$ws.Range("G7").Value = 18.48550382854987
$ws.Range("H7").Value = 11.83050654431597
$ws.Range("I7").Value = 16.53649501643003
$ws.Range("J7").Value = 7.307936951439714
$ws.Range("M7").Value = 18.84720710670399
$ws.Range("O7").Value = 16.61464603120016
$ws.Range("B8").Value = 11.83135035864673
$ws.Range("C8").Value = 9.253204952030735
$ws.Range("E8").Value = 23.25141775897441
$ws.Range("F8").Value = 35.99669934474885
$ws.Range("G8").Value = 18.31005381041722
$ws.Range("H8").Value = 11.75019884526946
$ws.Range("I8").Value = 16.35348582221637
$ws.Range("J8").Value = 7.289210175933114
$ws.Range("M8").Value = 19.21242394260137
$ws.Range("O8").Value = 16.4720034087217
$ws.Range("B9").Value = 13.48949190502424
$ws.Range("C9").Value = 10.20785983860373
$ws.Range("E9").Value = 23.34843964711696
$ws.Range("F9").Value = 36.04079492487867
$ws.Range("G9").Value = 18.0497446944023
$ws.Range("H9").Value = 11.61153593834539
$ws.Range("I9").Value = 16.02887280777987
$ws.Range("J9").Value = 7.25664142868387
$ws.Range("M9").Value = 19.92317005409375
$ws.Range("O9").Value = 16.23168535318849
$ws.Range("B10").Value = 14.5791195121037
$ws.Range("C10").Value = 10.85100872393255
$ws.Range("E10").Value = 23.44240639995836
$ws.Range("F10").Value = 36.12576944905005
$ws.Range("G10").Value = 17.91113474986555
$ws.Range("H10").Value = 11.52113110542411
$ws.Range("I10").Value = 15.81109735533567
$ws.Range("J10").Value = 7.235224515857869
$ws.Range("M10").Value = 20.43768193940945
$ws.Range("O10").Value = 16.07938023751209
$ws.Range("B11").Value = 15.04676010692907
$ws.Range("C11").Value = 11.13021842398306
$ws.Range("E11").Value = 23.48999571610281
$ws.Range("F11").Value = 36.17579700184368
$ws.Range("G11").Value = 17.85982512121676
$ws.Range("H11").Value = 11.48249450319102
$ws.Range("I11").Value = 15.71648793253985
$ws.Range("J11").Value = 7.226022223030435
$ws.Range("M11").Value = 20.66919590735545
$ws.Range("O11").Value = 16.01541195488235
$ws.Range("B12").Value = 15.21980742491598
$ws.Range("C12").Value = 11.23397828424744
$ws.Range("E12").Value = 23.50870441992573
$ws.Range("F12").Value = 36.19636889253367
$ws.Range("G12").Value = 17.84210865442279
$ws.Range("H12").Value = 11.4682219040308
$ws.Range("I12").Value = 15.68130019140627
$ws.Range("J12").Value = 7.222614917676474
$ws.Range("M12").Value = 20.75642831546692
$ws.Range("O12").Value = 15.99195782659782
$ws.Range("B13").Value = 15.18271823482475
$ws.Range("C13").Value = 11.21172009218949
$ws.Range("E13").Value = 23.50464473105883
$ws.Range("F13").Value = 36.19186612945641
$ws.Range("G13").Value = 17.84584763519674
$ws.Range("H13").Value = 11.47127982505825
$ws.Range("I13").Value = 15.68885012537996
$ws.Range("J13").Value = 7.223345304079396
$ws.Range("M13").Value = 20.73766179721795
$ws.Range("O13").Value = 15.99697480144195
$ws.Range("B14").Value = 15.06107776247423
$ws.Range("C14").Value = 11.13879456123492
$ws.Range("E14").Value = 23.4915211640985
$ws.Range("F14").Value = 36.17745689760315
$ws.Range("G14").Value = 17.85833310309272
$ws.Range("H14").Value = 11.48131310613246
$ws.Range("I14").Value = 15.71358022607205
$ws.Range("J14").Value = 7.225740352405801
$ws.Range("M14").Value = 20.6763817264812
$ws.Range("O14").Value = 16.01346692433085
$ws.Range("B15").Value = 14.98604347896989
$ws.Range("C15").Value = 11.09386758869244
$ws.Range("E15").Value = 23.4835718785192
$ws.Range("F15").Value = 36.16884251302664
$ws.Range("G15").Value = 17.86620464708167
$ws.Range("H15").Value = 11.4875054447035
$ws.Range("I15").Value = 15.72881125152748
$ws.Range("J15").Value = 7.22721745996492
$ws.Range("M15").Value = 20.63878693817096
$ws.Range("O15").Value = 16.02366915240858
$ws.Range("B16").Value = 14.54799327749438
$ws.Range("C16").Value = 10.83248842501632
$ws.Range("E16").Value = 23.43939305367891
$ws.Range("F16").Value = 36.12272825500223
$ws.Range("G16").Value = 17.91472641563005
$ws.Range("H16").Value = 11.52370621061804
$ws.Range("I16").Value = 15.81736981345722
$ws.Range("J16").Value = 7.235836755116116
$ws.Range("M16").Value = 20.42249465489047
$ws.Range("O16").Value = 16.08366809895838
$ws.Range("B17").Value = 14.27207675512299
$ws.Range("C17").Value = 10.66868062614685
$ws.Range("E17").Value = 23.41352528260452
$ws.Range("F17").Value = 36.09734651996864
$ws.Range("G17").Value = 17.94751893558348
$ws.Range("H17").Value = 11.54655192439082
$ws.Range("I17").Value = 15.87283775738665
$ws.Range("J17").Value = 7.241262599545902
$ws.Range("M17").Value = 20.28910390136931
$ws.Range("O17").Value = 16.12184049171283
$ws.Range("B18").Value = 14.11073798900615
$ws.Range("C18").Value = 10.5732073956872
$ws.Range("E18").Value = 23.39910328742958
$ws.Range("F18").Value = 36.08381862901884
$ws.Range("G18").Value = 17.96748441524342
$ws.Range("H18").Value = 11.55992643025292
$ws.Range("I18").Value = 15.90516116346166
$ws.Range("J18").Value = 7.244434281553344
$ws.Range("M18").Value = 20.21214609824396
$ws.Range("O18").Value = 16.14429632280164
$ws.Range("B19").Value = 14.05565849865145
$ws.Range("C19").Value = 10.54066784127755
$ws.Range("E19").Value = 23.39429890000083
$ws.Range("F19").Value = 36.07942246559171
$ws.Range("G19").Value = 17.97443324234078
$ws.Range("H19").Value = 11.56449503238402
$ws.Range("I19").Value = 15.91617746773882
$ws.Range("J19").Value = 7.245516906469509
$ws.Range("M19").Value = 20.18605138157895
$ws.Range("O19").Value = 16.1519852250042
$ws.Range("B20").Value = 14.30172171531134
$ws.Range("C20").Value = 10.68624855062548
$ws.Range("E20").Value = 23.41623176991982
$ws.Range("F20").Value = 36.09993765256256
$ws.Range("G20").Value = 17.94391366469279
$ws.Range("H20").Value = 11.5440957123248
$ws.Range("I20").Value = 15.86688967768886
$ws.Range("J20").Value = 7.240679745462939
$ws.Range("M20").Value = 20.30332842131819
$ws.Range("O20").Value = 16.11772518418546
$ws.Range("B21").Value = 15.09691613710391
$ws.Range("C21").Value = 11.16026839673493
$ws.Range("E21").Value = 23.49535728525168
$ws.Range("F21").Value = 36.18164514004837
$ws.Range("G21").Value = 17.85461912998875
$ws.Range("H21").Value = 11.47835636285802
$ws.Range("I21").Value = 15.70629907681465
$ws.Range("J21").Value = 7.2250347703375
$ws.Range("M21").Value = 20.69439358293953
$ws.Range("O21").Value = 16.0086018713589
$ws.Range("B22").Value = 15.59308502759218
$ws.Range("C22").Value = 11.45856316833509
$ws.Range("E22").Value = 23.5510739310241
$ws.Range("F22").Value = 36.24452711934163
$ws.Range("G22").Value = 17.80625715694342
$ws.Range("H22").Value = 11.43748019286905
$ws.Range("I22").Value = 15.605066218982
$ws.Range("J22").Value = 7.215260902480206
$ws.Range("M22").Value = 20.94739948807536
$ws.Range("O22").Value = 15.94176968297041
$ws.Range("B23").Value = 15.33042493178443
$ws.Range("C23").Value = 11.3004245336141
$ws.Range("E23").Value = 23.52097371067316
$ws.Range("F23").Value = 36.21010140626915
$ws.Range("G23").Value = 17.83114640190173
$ws.Range("H23").Value = 11.45910537754291
$ws.Range("I23").Value = 15.65875622369347
$ws.Range("J23").Value = 7.220436229273293
$ws.Range("M23").Value = 20.81262425785341
$ws.Range("O23").Value = 15.97702717539677
$ws.Range("B24").Value = 14.28832767149089
$ws.Range("C24").Value = 10.6783101309834
$ws.Range("E24").Value = 23.41500676544291
$ws.Range("F24").Value = 36.09876288537495
$ws.Range("G24").Value = 17.94554014379161
$ws.Range("H24").Value = 11.54520541760135
$ws.Range("I24").Value = 15.86957745198962
$ws.Range("J24").Value = 7.240943090880419
$ws.Range("M24").Value = 20.29689835035269
$ws.Range("O24").Value = 16.1195841264031
$ws.Range("B25").Value = 13.06343800502107
$ws.Range("C25").Value = 9.959544686137756
$ws.Range("E25").Value = 23.31818290058751
$ws.Range("F25").Value = 36.01962360480487
$ws.Range("G25").Value = 18.11103213131339
$ws.Range("H25").Value = 11.64703348110418
$ws.Range("I25").Value = 16.11303866356057
$ws.Range("J25").Value = 7.265009593096779
$ws.Range("M25").Value = 19.73193375801584
$ws.Range("O25").Value = 16.29245537891904
